# Results.docx edit script
# Applies the changes described by the commit's XML diff:
#  - "Studies included" count 8 -> 22
#  - Header / CI / I^2 lines: merge split runs (proofErr artifacts) back into single runs
#  - Forest-plot table: 8 study rows -> 22 study rows with updated values
#  - "Overall, IV" summary row values updated
#  - "Test of overall effect" line updated and surrounding blank paragraphs collapsed
#  - "Cochran's Q" / "H" rows updated; I^2 (%) row updated
#  - Final "I^2 = ..." paragraph split into two paragraphs (a line-wrap artifact)

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Phase A: plain text replacements addressed by the ORIGINAL
# (pre-edit) paragraph index. None of these change the paragraph
# count, so the indices stay valid for the whole phase.
# ---------------------------------------------------------------

# NOTE: a handful of paragraphs have more than one run (left behind
# by Word's grammar-checker proofErr splits). Assigning directly to
# a multi-run paragraph's own .Range.Text only overwrites the first
# run, so for those we rebuild an explicit Range(start,end) over the
# whole paragraph (its Start/End already span every run) before
# assigning - this both updates the text AND merges the runs back
# into one, matching the target XML.

$d.Paragraphs.Item(1).Range.Text = "Studies included: 22"

$p8 = $d.Paragraphs.Item(8)
$d.Range($p8.Range.Start, $p8.Range.End).Text = "Study                |     CFR     [95% Conf. Interval]   % Weight"

$d.Paragraphs.Item(10).Range.Text = "1                    |     0.250      0.200     1.250       0.43"

$d.Paragraphs.Item(19).Range.Text = "Overall, IV          |     1.989      1.872     2.113     100.00"

$d.Paragraphs.Item(21).Range.Text = "Test of overall effect = 1:  z =  22.314  p = 0.000"

$d.Paragraphs.Item(30).Range.Text = "Cochran's Q          |    703.16       21      0.000"

$p31 = $d.Paragraphs.Item(31)
$d.Range($p31.Range.Start, $p31.Range.End).Text = "                     |            -[95% Conf. Interval]-"

$d.Paragraphs.Item(32).Range.Text = "H                    |     5.787     5.358     6.207"

$p33 = $d.Paragraphs.Item(33)
$d.Range($p33.Range.Start, $p33.Range.End).Text = "I² (%)               |     97.0%     96.5%     97.4%"

$d.Paragraphs.Item(36).Range.Text = "I² = proportion of total variation in effect estimate due to between-study heterogene"

# ---------------------------------------------------------------
# Phase B: structural edits (paragraphs inserted / removed).
# Processed from the highest original index down to the lowest so
# that indices referenced later in this phase are never shifted by
# an edit performed earlier in the phase.
# ---------------------------------------------------------------

# --- split paragraph 36 into two paragraphs ---------------------
$p36 = $d.Paragraphs.Item(36)
$p36.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$d.Paragraphs.Item(37).Range.Text = "> ity (based on Q)"

# --- collapse the old "Test of overall effect" block ------------
# paragraph 21 already holds the new text (Phase A); paragraphs
# 22-24 (old text run + two blank paragraphs) are removed entirely.
$d = $word.ActiveDocument
$delRange = $d.Range($d.Paragraphs.Item(22).Range.Start, $d.Paragraphs.Item(24).Range.End)
$delRange.Delete()

# --- replace the 8-row study table with the 22-row version ------
# paragraph 10 already holds the updated row "1" (Phase A); delete
# old rows 2-8 (paragraphs 11-17) and insert the 21 new rows after.
$d = $word.ActiveDocument
$delRange2 = $d.Range($d.Paragraphs.Item(11).Range.Start, $d.Paragraphs.Item(17).Range.End)
$delRange2.Delete()

$newRows = @(
  "2                    |     0.200      0.100     0.300       1.21",
  "3                    |     2.150      1.400     6.600       0.61",
  "4                    |     1.700      1.240     3.180       1.65",
  "5                    |     0.800      0.570     9.400       0.19",
  "6                    |     0.900      0.100     3.200       0.12",
  "7                    |     2.600      1.500     3.900       1.60",
  "8                    |     1.800      0.700     3.900       0.49",
  "9                    |     1.900      1.500     2.370       6.97",
  "10                   |     1.700      1.040     5.200       0.56",
  "11                   |     0.700      0.580     0.830      11.36",
  "12                   |     2.700      1.700     4.000       1.99",
  "13                   |     4.000      2.170     5.250       1.87",
  "14                   |     2.100      1.440     2.450       5.17",
  "15                   |     3.600      2.400     5.200       2.44",
  "16                   |     3.200      1.200     6.900       0.48",
  "17                   |     4.700      4.130     5.350      21.78",
  "18                   |     5.400      4.300     6.700       7.42",
  "19                   |     0.600      0.500     0.700      12.89",
  "20                   |     4.100      2.800     6.000       2.51",
  "21                   |     1.000      0.240     1.100       0.63",
  "22                   |     2.000      1.800     2.400      17.63"
)

$d = $word.ActiveDocument
$anchor = $d.Paragraphs.Item(10)
foreach ($row in $newRows) {
    $anchor.Range.InsertParagraphAfter()
    $d = $word.ActiveDocument
    $anchor = $d.Paragraphs.Item($anchor.Index + 1)
    $anchor.Range.Text = $row
}

Write-Output "edit complete"
